$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Top5" column headers (new column O) for each of the 4 blocks ---
$ws.Range("O13").Value = "Top5"
$ws.Range("O24").Value = "Top5"
$ws.Range("O34").Value = "Top5"
$ws.Range("O43").Value = "Top5"

# --- Block 1 (E = 2): rows 14 (single) and 15:22 (fill range) ---
$ws.Range("L14").Formula = '=SMALL($A14:$I14, 2)'
$ws.Range("M14").Formula = '=SMALL($A14:$I14, 3)'
$ws.Range("N14").Formula = '=SMALL($A14:$I14, 4)'
$ws.Range("O14").Formula = '=SMALL($A14:$I14, 5)'

$ws.Range("L15:L22").Formula = '=SMALL($A15:$I15, 2)'
$ws.Range("M15:M22").Formula = '=SMALL($A15:$I15, 3)'
$ws.Range("N15:N22").Formula = '=SMALL($A15:$I15, 4)'
$ws.Range("O15:O22").Formula = '=SMALL($A15:$I15, 5)'

# --- Block 2 (E = 3): rows 25 (single) and 26:32 (fill range) ---
$ws.Range("L25").Formula = '=SMALL($A25:$I25, 2)'
$ws.Range("M25").Formula = '=SMALL($A25:$I25, 3)'
$ws.Range("N25").Formula = '=SMALL($A25:$I25, 4)'
$ws.Range("O25").Formula = '=SMALL($A25:$I25, 5)'

$ws.Range("L26:L32").Formula = '=SMALL($A26:$I26, 2)'
$ws.Range("M26:M32").Formula = '=SMALL($A26:$I26, 3)'
$ws.Range("N26:N32").Formula = '=SMALL($A26:$I26, 4)'
$ws.Range("O26:O32").Formula = '=SMALL($A26:$I26, 5)'

# --- Block 3 (E = 4): rows 35 (single) and 36:41 (fill range) ---
$ws.Range("L35").Formula = '=SMALL($A35:$I35, 2)'
$ws.Range("M35").Formula = '=SMALL($A35:$I35, 3)'
$ws.Range("N35").Formula = '=SMALL($A35:$I35, 4)'
$ws.Range("O35").Formula = '=SMALL($A35:$I35, 5)'

$ws.Range("L36:L41").Formula = '=SMALL($A36:$I36, 2)'
$ws.Range("M36:M41").Formula = '=SMALL($A36:$I36, 3)'
$ws.Range("N36:N41").Formula = '=SMALL($A36:$I36, 4)'
$ws.Range("O36:O41").Formula = '=SMALL($A36:$I36, 5)'

# --- Block 4 (E = 5): rows 44 (single) and 45:49 (fill range) ---
$ws.Range("L44").Formula = '=SMALL($A44:$I44, 2)'
$ws.Range("M44").Formula = '=SMALL($A44:$I44, 3)'
$ws.Range("N44").Formula = '=SMALL($A44:$I44, 4)'
$ws.Range("O44").Formula = '=SMALL($A44:$I44, 5)'

$ws.Range("L45:L49").Formula = '=SMALL($A45:$I45, 2)'
$ws.Range("M45:M49").Formula = '=SMALL($A45:$I45, 3)'
$ws.Range("N45:N49").Formula = '=SMALL($A45:$I45, 4)'
$ws.Range("O45:O49").Formula = '=SMALL($A45:$I45, 5)'

# --- Update the selection shown in the sheet view ---
$ws.Range("N54").Select()
